# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# F2: 261 -> 262
# F3: 377 -> 378

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 262
    $ws.Range("F3").Value = 378
}
